$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Stone" item/skill entry is inserted above the existing "aim" skill
# block (which currently occupies rows 35-36), so make room for it by
# inserting 3 rows at row 35: two rows for the new "stone_skill_title" /
# "stone_skill_description" entries plus one blank spacer row, matching
# the blank-row separator pattern used between every other block on this
# sheet (e.g. row 34, row 37, row 40 in the final layout).
$ws.Rows("35:37").Insert()

# Row 35: stone_skill_title
$ws.Range("A35").Value = "stone_skill_title"
# Row 36: stone_skill_description
$ws.Range("A36").Value = "stone_skill_description"

$ws.Range("B35").Value = "Stone"
$ws.Range("D35").Value = "Akmenis"
$ws.Range("C35").Value = "Камінчик"

$ws.Range("B36").Value = "Stones are usefull for building or throwing it. Please try not to throw all of them they are more usefull for small buildings."
$ws.Range("C36").Value = "Камінці корисні для будівництва або кидання. Будь ласка, намагайтеся не викидати їх усі, вони більш корисні для невеликих будівель."
$ws.Range("D36").Value = "Akmenys naudingi statant ar mėtant. Stenkitės neišmesti jų visų, jie yra naudingesni mažiems pastatams."

# The description row wraps onto multiple lines, same as the other
# description rows on this sheet, so it needs a taller row height.
$ws.Rows(36).RowHeight = 60

# Reflect where the author was last working when they saved.
$ws.Range("E36").Select()
